$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell H1, matching style of existing header cells (A1:G1 use style index 1)
$ws.Range("H1").Value = "Koppel blijft bijelkaar"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# Values for H2:H111 ("Koppel blijft bijelkaar" flag), one integer per data row
$values = @(1,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,1,0,0,0,0,0,0,0,0,0,0,0,1,1,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0)

if ($values.Length -ne 110) {
    throw "Expected 110 values, got $($values.Length)"
}

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $values[$i]
}
